$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two data rows (original rows 2 and 3) turned out to be a
# "false start" - remove them. Select the two row headers (as a user
# would) and delete them outright; the rows below (old rows 4 and 5,
# the real data) shift up to become the new rows 2 and 3.
$ws.Rows("2:3").Select() | Out-Null
$excel.Selection.EntireRow.Delete() | Out-Null
